# Refresh the cryptocurrency price/volume snapshot in-place.
# Coin name/link cells are only touched where coins swapped rank position;
# price (D) and volume (E) cells are refreshed with the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column holds formatted text (e.g. "58.667.17"); force Text format
# first so Excel does not reinterpret numeric-looking values as numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "58.667.17"
$ws.Range("E2").Value = "  +0.82%  "

# Row 3
$ws.Range("D3").Value = "3.154.57"
$ws.Range("E3").Value = "  +0.48%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "529.29"
$ws.Range("E5").Value = "  -0.95%  "

# Row 6
$ws.Range("D6").Value = "139.38"
$ws.Range("E6").Value = "  +0.18%  "

# Row 7
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.11%  "

# Row 8
$ws.Range("D8").Value = "0.536"
$ws.Range("E8").Value = "  +14.45%  "

# Row 9
$ws.Range("D9").Value = "7.31"
$ws.Range("E9").Value = "  +0.17%  "

# Row 10
$ws.Range("D10").Value = "0.438"
$ws.Range("E10").Value = "  +5.66%  "

# Row 11
$ws.Range("E11").Value = "  +3.87%  "

# Row 12
$ws.Range("E12").Value = "  +3.21%  "

# Row 13
$ws.Range("D13").Value = "3.699.72"

# Row 14
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "0.0000172"
$ws.Range("E14").Value = "  +4.46%  "

# Row 15
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "25.74"
$ws.Range("E15").Value = "  +0.87%  "

# Row 16
$ws.Range("D16").Value = "58.715.88"
$ws.Range("E16").Value = "  +0.72%  "

# Row 17
$ws.Range("D17").Value = "6.25"
$ws.Range("E17").Value = "  +3.26%  "

# Row 18
$ws.Range("D18").Value = "3.148.57"
$ws.Range("E18").Value = "  +0.25%  "

# Row 19
$ws.Range("D19").Value = "12.97"
$ws.Range("E19").Value = "  +2.11%  "

# Row 20
$ws.Range("D20").Value = "8.11"
$ws.Range("E20").Value = "  -0.61%  "

# Row 21
$ws.Range("D21").Value = "371.89"
$ws.Range("E21").Value = "  +3.53%  "

# Row 22
$ws.Range("E22").Value = "  +1.66%  "

# Row 23
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.17%  "

# Row 24
$ws.Range("D24").Value = "0.531"
$ws.Range("E24").Value = "  +4.85%  "

# Row 25
$ws.Range("D25").Value = "69.56"
$ws.Range("E25").Value = "  +0.83%  "

# Row 26
$ws.Range("D26").Value = "0.167"
$ws.Range("E26").Value = "  +0.39%  "

# Row 27
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.34%  "

# Row 28
$ws.Range("D28").Value = "8.30"
$ws.Range("E28").Value = "  +13.45%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0863"
$ws.Range("E29").Value = "  -1.72%  "

# Row 30
$ws.Range("D30").Value = "22.20"
$ws.Range("E30").Value = "  +2.63%  "

# Row 31
$ws.Range("E31").Value = "  -0.15%  "

# Row 32
$ws.Range("D32").Value = "6.07"
$ws.Range("E32").Value = "  -0.87%  "

# Row 33
$ws.Range("D33").Value = "5.13"
$ws.Range("E33").Value = "  +1.28%  "

# Row 34
$ws.Range("E34").Value = "  +0.67%  "

# Row 35
$ws.Range("D35").Value = "6.29"
$ws.Range("E35").Value = "  +3.34%  "

# Row 36
$ws.Range("D36").Value = "158.10"
$ws.Range("E36").Value = "  -0.75%  "

# Row 37
$ws.Range("E37").Value = "  +5.25%  "

# Row 38
$ws.Range("D38").Value = "25.06"
$ws.Range("E38").Value = "  -2.94%  "

# Row 39
$ws.Range("D39").Value = "1.68"
$ws.Range("E39").Value = "  -0.60%  "

# Row 40
$ws.Range("D40").Value = "0.0685"
$ws.Range("E40").Value = "  +1.96%  "

# Row 41
$ws.Range("D41").Value = "2.617.22"
$ws.Range("E41").Value = "  +5.16%  "

# Row 42
$ws.Range("D42").Value = "4.23"
$ws.Range("E42").Value = "  +5.65%  "

# Row 43
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "0.724"
$ws.Range("E43").Value = "  +3.01%  "

# Row 44
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "39.06"
$ws.Range("E44").Value = "  +4.28%  "

# Row 45
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "0.0286"
$ws.Range("E45").Value = "  +6.69%  "

# Row 46
$ws.Range("E46").Value = "  -0.09%  "

# Row 47
$ws.Range("D47").Value = "3.195.31"
$ws.Range("E47").Value = "  +0.42%  "

# Row 48
$ws.Range("D48").Value = "0.103"
$ws.Range("E48").Value = "  +12.78%  "

# Row 49
$ws.Range("E49").Value = "  +1.94%  "

# Row 50
$ws.Range("D50").Value = "0.978"
$ws.Range("E50").Value = "  -1.70%  "

# Row 51
$ws.Range("D51").Value = "20.12"
$ws.Range("E51").Value = "  +1.34%  "
